$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "287.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.73%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.12"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.200"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.10%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06958"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.12%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.430"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.32%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.09%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.397"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.59%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8993"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.96%"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.66%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07636"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "15.33%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07736"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.98%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02933"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.62%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09012"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.33%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001574"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.10%"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "One"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006497"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.94%"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006420"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.11%"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.490"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.14%"

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.230"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.08%"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3233"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.87%"

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1337"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.38%"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.042"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.44%"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1599"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.00%"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04532"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.36%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001210"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.64%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.46%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001170"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.06%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001667"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "3.40%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04353"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.66%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006938"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.45%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.22%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002071"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.85%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01172"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.80%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005819"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.61%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.20%"
